$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.101.54'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -0.99%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.943.70'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -1.59%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '375.05'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -1.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.29'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -3.56%  '

$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("E9").Value = '  -2.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.35'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -3.07%  '

$ws.Range("E11").Value = '  -0.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0850'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.416.62'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -1.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.06'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -2.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.58'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.945.59'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.996'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +1.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.70'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +43.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.056.48'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -0.98%  '

$ws.Range("E20").Value = '  -7.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.43'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -4.49%  '

$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '265.66'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.69'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -1.10%  '

$ws.Range("E25").Value = '  +7.93%  '

$ws.Range("E26").Value = '  -2.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.58'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -2.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("B29").Value = 'EthereumClassic'

$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.62'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -1.77%  '

$ws.Range("B30").Value = 'Kaspa'

$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.164'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -5.05%  '

$ws.Range("E31").Value = '  -5.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.03'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +1.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.75'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -0.60%  '

$ws.Range("E34").Value = '  -1.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.35'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -5.25%  '

$ws.Range("E36").Value = '  -3.14%  '

$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("E38").Value = '  +3.00%  '

$ws.Range("E39").Value = '  -1.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.37'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -5.09%  '

$ws.Range("E41").Value = '  -3.43%  '

$ws.Range("E42").Value = '  -4.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.40'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -4.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.35'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -2.38%  '

$ws.Range("E45").Value = '  -0.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.33'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.272'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -4.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.33'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -1.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.990.35'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -2.61%  '

$ws.Range("E50").Value = '  -3.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.32'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +1.42%  '
